# DPLKINV081-001 - update "No. Urut" value from 1369 to 2962
# and adjust the active view (scrolled to column Q, selection on AB2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the numeric "No. Urut" cell (M2) to the new value.
$ws.Range("M2").Value2 = 2962

# Update the big PREPARATION text block (F2) so the embedded
# "No. Urut : 1369" becomes "No. Urut : 2962".
$prepText = $ws.Range("F2").Value2
$prepText = $prepText -replace "No\. Urut : 1369", "No. Urut : 2962"
$ws.Range("F2").Value2 = $prepText

# Update the sheet view: scroll so column Q is the left-most visible
# column, and select cell AB2.
$ws.Application.ActiveWindow.ScrollColumn = 17
$ws.Range("AB2").Select()
